$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the values from column D (rows 2-13) into column E (rows 2-13)
for ($row = 2; $row -le 13; $row++) {
    $dValue = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 5).Value = $dValue
}

# Update the active selection to E14
$ws.Range("E14").Select()
